$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLine = [char]10
$newText = "1.广播域限缩在源/目的节点内的矩形内，当长距离传输时，广播域会非常大，" + $newLine + "2.考虑到排队时延，局部网络拥塞可能导致广播域内找不到时延低的路径" + $newLine

$ws.Range("E6").Value = $newText
